$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for all existing data rows
# (rows 2 through 295) from 45205 to 45206.
$ws.Range("C2:C295").Value = 45206

# Force row 295 to carry an explicit row height (matches diff: ht="15" customHeight="1").
$ws.Rows.Item(295).RowHeight = 15

# Append the new record as row 296.
$newRow = 296
$ws.Cells.Item($newRow, 1).Value = "A 48035-2023"

$ws.Cells.Item($newRow, 2).Value = 45204
$ws.Cells.Item($newRow, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($newRow, 3).Value = 45206
$ws.Cells.Item($newRow, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($newRow, 4).Value = "VÄSTRA GÖTALANDS LÄN"
$ws.Cells.Item($newRow, 5).Value = "TIDAHOLM"

$ws.Cells.Item($newRow, 7).Value = 1.4
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 0
$ws.Cells.Item($newRow, 11).Value = 0
$ws.Cells.Item($newRow, 12).Value = 0
$ws.Cells.Item($newRow, 13).Value = 0
$ws.Cells.Item($newRow, 14).Value = 0
$ws.Cells.Item($newRow, 15).Value = 0
$ws.Cells.Item($newRow, 16).Value = 0
$ws.Cells.Item($newRow, 17).Value = 0

# Column R keeps the wrap-text style used by the rest of the table, but stays empty.
$ws.Cells.Item($newRow, 18).WrapText = $true
